$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 0.000346637487411499
$ws.Range("M2").Value = 0.07063746452331543
$ws.Range("N2").Value = 0.0001251697540283203

$ws.Range("L3").Value = 0.0000729827880859375
$ws.Range("M3").Value = 0.0004432201385498047
$ws.Range("N3").Value = 0.00005936622619628906

$ws.Range("E4").Value = 0.8502994011976048
$ws.Range("F4").Value = 0.3125
$ws.Range("G4").Value = 0.2631578947368421
$ws.Range("H4").Value = 137.0
$ws.Range("I4").Value = 11.0
$ws.Range("J4").Value = 14.0
$ws.Range("K4").Value = 5.0
$ws.Range("L4").Value = 0.00029168701171875
$ws.Range("M4").Value = 0.001210212707519531
$ws.Range("N4").Value = 0.0001287460327148438

$ws.Range("E5").Value = 0.8622754491017964
$ws.Range("F5").Value = 0.375
$ws.Range("G5").Value = 0.3157894736842105
$ws.Range("H5").Value = 138.0
$ws.Range("I5").Value = 10.0
$ws.Range("J5").Value = 13.0
$ws.Range("K5").Value = 6.0
$ws.Range("L5").Value = 0.00009527349472045898
$ws.Range("M5").Value = 0.0005195140838623047
$ws.Range("N5").Value = 0.00006747245788574219

$ws.Range("E6").Value = 0.8862275449101796
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.3684210526315789
$ws.Range("J6").Value = 12.0
$ws.Range("K6").Value = 7.0
$ws.Range("L6").Value = 0.01129698610305786
$ws.Range("M6").Value = 0.0662240982055664
$ws.Range("N6").Value = 0.007999181747436523

$ws.Range("E7").Value = 0.8982035928143712
$ws.Range("F7").Value = 0.5625
$ws.Range("H7").Value = 141.0
$ws.Range("I7").Value = 7.0
$ws.Range("L7").Value = 0.007614052534103394
$ws.Range("M7").Value = 0.013092041015625
$ws.Range("N7").Value = 0.006600379943847656

$ws.Range("L8").Value = 0.0003924052715301514
$ws.Range("M8").Value = 0.002433061599731445
$ws.Range("N8").Value = 0.0001382827758789062

$ws.Range("L9").Value = 0.0001002638339996338
$ws.Range("M9").Value = 0.0006198883056640625
$ws.Range("N9").Value = 0.00008082389831542969

$ws.Range("L10").Value = 0.0004473221302032471
$ws.Range("M10").Value = 0.01003742218017578
$ws.Range("N10").Value = 0.0001530647277832031

$ws.Range("L11").Value = 0.0001359720230102539
$ws.Range("M11").Value = 0.0006124973297119141
$ws.Range("N11").Value = 0.0001018047332763672

$ws.Range("L12").Value = 0.001205533266067505
$ws.Range("M12").Value = 0.004215717315673828
$ws.Range("N12").Value = 0.00067138671875

$ws.Range("L13").Value = 0.0005541503429412842
$ws.Range("M13").Value = 0.001319408416748047
$ws.Range("N13").Value = 0.0004651546478271484

$ws.Range("L14").Value = 0.000302384614944458
$ws.Range("M14").Value = 0.001739740371704102
$ws.Range("N14").Value = 0.0001285076141357422

$ws.Range("L15").Value = 0.00008963632583618164
$ws.Range("M15").Value = 0.0005412101745605469
$ws.Range("N15").Value = 0.00006198883056640625

$ws.Range("L16").Value = 0.000583186149597168
$ws.Range("M16").Value = 0.006662130355834961
$ws.Range("N16").Value = 0.0002110004425048828

$ws.Range("L17").Value = 0.0001706204414367676
$ws.Range("M17").Value = 0.002483606338500977
$ws.Range("N17").Value = 0.0001039505004882812
